# "Extracted data from excel spreadsheet"
#
# Cell C4 on the "TestData" sheet used to hold the text "asdfaghh" (a
# shared string). The data-extraction step replaced it with the numeric
# value 15. Because that string is then unreferenced by any remaining
# cell, it drops out of the shared-strings table on save and every
# entry that followed it re-indexes down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 15

# Leave the final active selection on C5, matching the saved view state.
$null = $ws.Range("C5").Select()
